$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.901.54"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -0.48%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'1.804.79"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -0.88%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'1.000"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'  -0.14%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'309.88"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -0.46%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'0.9999"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  -0.10%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'0.4410"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  +4.47%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("E8").Value = "'  +0.89%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.07466"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  +3.54%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.8596"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  +2.35%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("E11").Value = "'  -0.17%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'1.801.78"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  -1.04%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'6.639"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  -0.33%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'93.17"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  +3.82%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'0.07068"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  +0.08%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'5.277"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  -0.10%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("E17").Value = "'  -0.24%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'0.000008702"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  -0.28%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("E19").Value = "'  -0.11%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'14.82"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  -0.44%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'26.931.88"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  -0.60%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'5.165"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  +0.52%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'10.82"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  +0.29%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("B24").Value = "'WrappedliquidstakedEther2.0"
$ws.Range("B24").Style = "Normal"
$ws.Range("C24").Value = "'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("C24").Style = "Normal"
$ws.Range("D24").Value = "'2.016.92"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  -1.68%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("B25").Value = "'Toncoin"
$ws.Range("B25").Style = "Normal"
$ws.Range("C25").Value = "'https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("C25").Style = "Normal"
$ws.Range("D25").Value = "'1.989"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  +0.54%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("B26").Value = "'Monero"
$ws.Range("B26").Style = "Normal"
$ws.Range("C26").Value = "'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("C26").Style = "Normal"
$ws.Range("D26").Value = "'151.26"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  -0.41%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("B27").Value = "'LidoDAOToken"
$ws.Range("B27").Style = "Normal"
$ws.Range("C27").Value = "'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("C27").Style = "Normal"
$ws.Range("D27").Value = "'2.208"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  -1.25%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("B28").Value = "'EthereumClassic"
$ws.Range("B28").Style = "Normal"
$ws.Range("C28").Value = "'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("C28").Style = "Normal"
$ws.Range("D28").Value = "'18.36"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  +0.75%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("B29").Value = "'InternetComputer(DFINITY)"
$ws.Range("B29").Style = "Normal"
$ws.Range("C29").Value = "'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("C29").Style = "Normal"
$ws.Range("D29").Value = "'5.203"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  -0.74%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("B30").Value = "'BitcoinCash"
$ws.Range("B30").Style = "Normal"
$ws.Range("C30").Value = "'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("C30").Style = "Normal"
$ws.Range("D30").Value = "'117.31"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  +0.49%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("B31").Value = "'Stellar"
$ws.Range("B31").Style = "Normal"
$ws.Range("C31").Value = "'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("C31").Style = "Normal"
$ws.Range("D31").Value = "'0.08778"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  +0.77%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("B32").Value = "'ImmutableX"
$ws.Range("B32").Style = "Normal"
$ws.Range("C32").Value = "'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("C32").Style = "Normal"
$ws.Range("D32").Value = "'0.7426"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  +1.21%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("B33").Value = "'ARBITRUM"
$ws.Range("B33").Style = "Normal"
$ws.Range("C33").Value = "'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("C33").Style = "Normal"
$ws.Range("D33").Value = "'1.164"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  -0.70%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("B34").Value = "'Filecoin"
$ws.Range("B34").Style = "Normal"
$ws.Range("C34").Value = "'https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("C34").Style = "Normal"
$ws.Range("D34").Value = "'4.476"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  +1.57%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("B35").Value = "'HuobiToken"
$ws.Range("B35").Style = "Normal"
$ws.Range("C35").Value = "'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("C35").Style = "Normal"
$ws.Range("D35").Value = "'2.889"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  +0.25%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("B36").Value = "'Frax"
$ws.Range("B36").Style = "Normal"
$ws.Range("C36").Value = "'https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("C36").Style = "Normal"
$ws.Range("D36").Value = "'0.9995"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  -0.15%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("B37").Value = "'TrustWalletToken"
$ws.Range("B37").Style = "Normal"
$ws.Range("C37").Value = "'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("C37").Style = "Normal"
$ws.Range("D37").Value = "'1.097"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  +0.51%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("B38").Value = "'VeChain"
$ws.Range("B38").Style = "Normal"
$ws.Range("C38").Value = "'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("C38").Style = "Normal"
$ws.Range("D38").Value = "'0.01972"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  +1.37%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("B39").Value = "'Hedera"
$ws.Range("B39").Style = "Normal"
$ws.Range("C39").Value = "'https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("C39").Style = "Normal"
$ws.Range("D39").Value = "'0.05207"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  -0.85%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("B40").Value = "'TheSandbox"
$ws.Range("B40").Style = "Normal"
$ws.Range("C40").Value = "'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("C40").Style = "Normal"
$ws.Range("D40").Value = "'0.5249"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  +4.12%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("B41").Value = "'FraxShare"
$ws.Range("B41").Style = "Normal"
$ws.Range("C41").Value = "'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("C41").Style = "Normal"
$ws.Range("D41").Value = "'7.067"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  -3.17%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("B42").Value = "'MXToken"
$ws.Range("B42").Style = "Normal"
$ws.Range("C42").Value = "'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("C42").Style = "Normal"
$ws.Range("D42").Value = "'2.820"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  -1.57%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("B43").Value = "'Algorand"
$ws.Range("B43").Style = "Normal"
$ws.Range("C43").Value = "'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("C43").Style = "Normal"
$ws.Range("D43").Value = "'0.1682"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  +0.01%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("B44").Value = "'Aptos"
$ws.Range("B44").Style = "Normal"
$ws.Range("C44").Value = "'https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("C44").Style = "Normal"
$ws.Range("D44").Value = "'8.492"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  -0.48%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("B45").Value = "'RenderToken"
$ws.Range("B45").Style = "Normal"
$ws.Range("C45").Value = "'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("C45").Style = "Normal"
$ws.Range("D45").Value = "'2.063"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  +7.50%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("B46").Value = "'Decentraland"
$ws.Range("B46").Style = "Normal"
$ws.Range("C46").Value = "'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("C46").Style = "Normal"
$ws.Range("D46").Value = "'0.4974"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  +5.78%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("B47").Value = "'EnergySwap"
$ws.Range("B47").Style = "Normal"
$ws.Range("C47").Value = "'https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("C47").Style = "Normal"
$ws.Range("D47").Value = "'10.45"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  -1.01%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("B48").Value = "'Quant"
$ws.Range("B48").Style = "Normal"
$ws.Range("C48").Value = "'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("C48").Style = "Normal"
$ws.Range("D48").Value = "'104.10"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  -1.74%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("B49").Value = "'PaxDollar"
$ws.Range("B49").Style = "Normal"
$ws.Range("C49").Value = "'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("C49").Style = "Normal"
$ws.Range("D49").Value = "'0.9992"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  -0.10%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("B50").Value = "'NEARProtocol"
$ws.Range("B50").Style = "Normal"
$ws.Range("C50").Value = "'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("C50").Style = "Normal"
$ws.Range("D50").Value = "'1.666"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  +1.24%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("B51").Value = "'Cronos"
$ws.Range("B51").Style = "Normal"
$ws.Range("C51").Value = "'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("C51").Style = "Normal"
$ws.Range("D51").Value = "'0.06346"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  +0.53%  "
$ws.Range("E51").Style = "Normal"
